{"js": "const body = context.document.body;\n\n// 1) First date near the top: \"07 \u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\" -> \"07 \u0436\u043e\u0432\u0442\u043d\u044f 2024 \u0440\u043e\u043a\u0443\"\nconst firstDate = body.search(\"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\", { matchCase: true, matchWholeWord: false });\nfirstDate.load(\"text\");\nawait context.sync();\nfirstDate.items[0].insertText(\"\u0436\u043e\u0432\u0442\u043d\u044f 2024 \u0440\u043e\u043a\u0443\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 2) \"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\" + [_GoBack bookmark] + \"\u043d\u044f \u0440\u043e\u0431\u043e\u0442\u0438:\" -> single run\n//    \"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\u043d\u044f \u0440\u043e\u0431\u043e\u0442\u0438:\" (bookmark removed).\nconst outputsLabel = body.search(\"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\", { matchCase: true, matchWholeWord: false });\noutputsLabel.load(\"paragraphs\");\nawait context.sync();\nconst outputsParagraph = outputsLabel.items[0].paragraphs.getFirst();\noutputsParagraph.getRange().insertText(\"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\u043d\u044f \u0440\u043e\u0431\u043e\u0442\u0438:\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 3) Issue date of the assignment: \"  01 \u0431\u0435\u0440\u0435\u0437\u043d\u044f 2023 \u0440\u043e\u043a\u0443\" -> \"  07 \u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\"\nconst issueDate = body.search(\"01 \u0431\u0435\u0440\u0435\u0437\u043d\u044f 2023 \u0440\u043e\u043a\u0443\", { matchCase: true, matchWholeWord: false });\nissueDate.load(\"text\");\nawait context.sync();\nissueDate.items[0].insertText(\"07 \u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\", Word.InsertLocation.replace);\nawait context.sync();\n\n// 4) Last date (end of document): \"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443 \" -> \"\u0436\u043e\u0432\u0442\u043d\u044f\" + [_GoBack bookmark] + \" 2024 \u0440\u043e\u043a\u0443 \"\nconst lastMonth = body.search(\"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430\", { matchCase: true, matchWholeWord: false });\nlastMonth.load(\"text\");\nawait context.sync();\nconst lastMonthRange = lastMonth.items[lastMonth.items.length - 1];\n\n// Put the (new) _GoBack bookmark where the cursor would land after typing,\n// i.e. right after \"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430\" and before \" 2024 \u0440\u043e\u043a\u0443 \", then rename the\n// month text in place so the surrounding run formatting is preserved.\nconst splitPoint = lastMonthRange.getRange(Word.RangeLocation.end);\nsplitPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n\nlastMonthRange.insertText(\"\u0436\u043e\u0432\u0442\u043d\u044f\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) First date near the top: \"07 \u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\" -> \"07 \u0436\u043e\u0432\u0442\u043d\u044f 2024 \u0440\u043e\u043a\u0443\"\n$find1 = $d.Content.Find\n$find1.Execute(\"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\", $false, $false, $false, $false, $false, $true, 1, $false, \"\u0436\u043e\u0432\u0442\u043d\u044f 2024 \u0440\u043e\u043a\u0443\", 1)\n\n# 2) \"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\" + [_GoBack bookmark] + \"\u043d\u044f \u0440\u043e\u0431\u043e\u0442\u0438:\" -> single run\n#    \"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\u043d\u044f \u0440\u043e\u0431\u043e\u0442\u0438:\" (the old _GoBack bookmark is removed as a\n#    side effect of Find/Replace spanning across it).\n$find2 = $d.Content.Find\n$find2.Execute(\"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\" + \"\u043d\u044f \u0440\u043e\u0431\u043e\u0442\u0438:\", $false, $false, $false, $false, $false, $true, 1, $false, \"4.\u0412\u0438\u0445\u0456\u0434\u043d\u0456  \u0434\u0430\u043d\u0456 \u0434\u043e \u0432\u0438\u043a\u043e\u043d\u0430\u043d\u043d\u044f \u0440\u043e\u0431\u043e\u0442\u0438:\", 2)\n\n# 3) Issue date of the assignment: \"  01 \u0431\u0435\u0440\u0435\u0437\u043d\u044f 2023 \u0440\u043e\u043a\u0443\" -> \"  07 \u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\"\n$find3 = $d.Content.Find\n$find3.Execute(\"01 \u0431\u0435\u0440\u0435\u0437\u043d\u044f 2023 \u0440\u043e\u043a\u0443\", $false, $false, $false, $false, $false, $true, 1, $false, \"07 \u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443\", 1)\n\n# 4) Last date (end of document): \"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430 2024 \u0440\u043e\u043a\u0443 \" -> \"\u0436\u043e\u0432\u0442\u043d\u044f\" + [_GoBack bookmark] + \" 2024 \u0440\u043e\u043a\u0443 \"\n#    Locate the LAST occurrence of \"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430\" (the one near the very end of the\n#    document) by scanning forward and remembering the final match.\n$fullEnd = $d.Content.End\n$lastStart = -1\n$lastEnd = -1\n$scanRange = $d.Range(0, $fullEnd)\n$scanFind = $scanRange.Find\n$more = $true\nwhile ($more) {\n  $scanFind.Execute(\"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430\")\n  if ($scanFind.Found) {\n    $lastStart = $scanRange.Start\n    $lastEnd = $scanRange.End\n    $scanRange.Start = $scanRange.End\n    $scanRange.End = $d.Content.End\n  } else {\n    $more = $false\n  }\n}\n\n$monthRange = $d.Range($lastStart, $lastEnd)\n\n# Put the (new) _GoBack bookmark where the cursor would land after typing, i.e.\n# right after \"\u043b\u0438\u0441\u0442\u043e\u043f\u0430\u0434\u0430\" and before \" 2024 \u0440\u043e\u043a\u0443 \".\n$splitPoint = $d.Range($lastEnd, $lastEnd)\n$d.Bookmarks.Add(\"_GoBack\", $splitPoint)\n\n# Rename the month text in place so the surrounding run formatting is preserved.\n$monthRange.Text = \"\u0436\u043e\u0432\u0442\u043d\u044f\"\n"}
